$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.380.72"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "3.602.69"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.15"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "190.13"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.632"
$ws.Range("E7").Value = "  -2.26%  "
$ws.Range("D8").Value = "3.599.88"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.185"
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.662"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.10"
$ws.Range("E12").Value = "  -3.85%  "
$ws.Range("E13").Value = "  +7.56%  "
$ws.Range("E14").Value = "  -2.08%  "
$ws.Range("D15").Value = "4.181.47"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.82"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "3.600.42"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "70.346.70"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.65"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "491.20"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.54"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("E24").Value = "  -9.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.75"
$ws.Range("E25").Value = "  +6.30%  "
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("E27").Value = "  -4.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.03"
$ws.Range("E28").Value = "  -2.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.38"
$ws.Range("E29").Value = "  -3.10%  "
$ws.Range("E30").Value = "  -2.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.62"
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.27"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "66.21"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.117"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "577.22"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.82"
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("D37").Value = "0.0₃0816"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("E39").Value = "  -3.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.96"
$ws.Range("E40").Value = "  +6.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.24"
$ws.Range("E41").Value = "  +17.79%  "
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("E43").Value = "  -6.39%  "
$ws.Range("D44").Value = "3.223.05"
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.41"
$ws.Range("E47").Value = "  +3.59%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.60"
$ws.Range("E48").Value = "  +4.10%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  -3.32%  "
